$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the original LEILA row (account 004382902, balance 26000)
# that currently sits at sheet row 3 (right after the header row and
# LEONARDO). Deleting it shifts THIAGO/LUIS/TATIANA/BLUEMETRIX/... up by
# one row.
$ws.Rows.Item(3).Delete()

# Step 2: after the shift, row 5 holds TATIANA (005348011) and row 6 holds
# BLUEMETRIX (001761119). Both accounts are being dropped from the export;
# overwrite those two rows in place with LEILA's updated balance and with
# PEDRO's updated balance (he used to appear further down the sheet).
# The leading apostrophe forces the account numbers to stay text so the
# leading zeros are preserved (they are account numbers, not numbers).
$ws.Range("A5").Value = "'004382902"
$ws.Range("B5").Value = "LEILA"
$ws.Range("C5").Value = 10000

$ws.Range("A6").Value = "'005324840"
$ws.Range("B6").Value = "PEDRO"
$ws.Range("C6").Value = 4734.48

# Step 3: delete PEDRO's old row (account 005324840, balance 734.48). It
# originally lived at sheet row 21, but the row-3 deletion in step 1
# already shifted everything below it up by one, so it is now at row 20.
$ws.Rows.Item(20).Delete()
